# Apply updated cryptocurrency price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '''21.709.38'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '''  +5.86%  '
$ws.Range('E2').Style = 'Normal'

# Row 3
$ws.Range('D3').Value = '''1.565.62'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '''  +6.38%  '
$ws.Range('E3').Style = 'Normal'

# Row 4
$ws.Range('D4').Value = '''1.002'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '''  -0.66%  '
$ws.Range('E4').Style = 'Normal'

# Row 5
$ws.Range('D5').Value = '''0.9829'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '''  +2.68%  '
$ws.Range('E5').Style = 'Normal'

# Row 6
$ws.Range('D6').Value = '''284.19'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '''  +2.72%  '
$ws.Range('E6').Style = 'Normal'

# Row 7
$ws.Range('D7').Value = '''0.3671'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '''  +0.68%  '
$ws.Range('E7').Style = 'Normal'

# Row 8
$ws.Range('D8').Value = '''0.3247'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '''  +6.27%  '
$ws.Range('E8').Style = 'Normal'

# Row 9
$ws.Range('D9').Value = '''41.21'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '''  +3.85%  '
$ws.Range('E9').Style = 'Normal'

# Row 10
$ws.Range('D10').Value = '''1.118'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '''  +6.08%  '
$ws.Range('E10').Style = 'Normal'

# Row 11
$ws.Range('D11').Value = '''0.07015'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '''  +6.11%  '
$ws.Range('E11').Style = 'Normal'

# Row 12
$ws.Range('D12').Value = '''0.9986'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '''  -0.34%  '
$ws.Range('E12').Style = 'Normal'

# Row 13
$ws.Range('D13').Value = '''19.91'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '''  +9.82%  '
$ws.Range('E13').Style = 'Normal'

# Row 14
$ws.Range('D14').Value = '''5.786'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '''  +6.04%  '
$ws.Range('E14').Style = 'Normal'

# Row 15
$ws.Range('D15').Value = '''6.459'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '''  +4.64%  '
$ws.Range('E15').Style = 'Normal'

# Row 16
$ws.Range('B16').Value = 'Dai'
$ws.Range('C16').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D16').Value = '''0.9834'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '''  +2.06%  '
$ws.Range('E16').Style = 'Normal'

# Row 17
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').Value = '''0.00001060'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '''  +3.14%  '
$ws.Range('E17').Style = 'Normal'

# Row 18
$ws.Range('D18').Value = '''1.560.50'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '''  +5.89%  '
$ws.Range('E18').Style = 'Normal'

# Row 19
$ws.Range('E19').Value = '''  +4.60%  '
$ws.Range('E19').Style = 'Normal'

# Row 20
$ws.Range('D20').Value = '''73.80'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '''  +6.94%  '
$ws.Range('E20').Style = 'Normal'

# Row 21
$ws.Range('D21').Value = '''16.02'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '''  +10.62%  '
$ws.Range('E21').Style = 'Normal'

# Row 22
$ws.Range('D22').Value = '''5.815'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '''  +6.60%  '
$ws.Range('E22').Style = 'Normal'

# Row 23
$ws.Range('D23').Value = '''11.53'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '''  +4.54%  '
$ws.Range('E23').Style = 'Normal'

# Row 24
$ws.Range('D24').Value = '''21.757.06'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '''  +5.81%  '
$ws.Range('E24').Style = 'Normal'

# Row 25
$ws.Range('D25').Value = '''2.350'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '''  +4.41%  '
$ws.Range('E25').Style = 'Normal'

# Row 26
$ws.Range('D26').Value = '''2.375'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '''  +11.93%  '
$ws.Range('E26').Style = 'Normal'

# Row 27
$ws.Range('D27').Value = '''148.16'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '''  +5.30%  '
$ws.Range('E27').Style = 'Normal'

# Row 28
$ws.Range('D28').Value = '''18.15'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '''  +5.56%  '
$ws.Range('E28').Style = 'Normal'

# Row 29
$ws.Range('D29').Value = '''1.734.98'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '''  +6.36%  '
$ws.Range('E29').Style = 'Normal'

# Row 30
$ws.Range('D30').Value = '''119.98'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '''  +5.55%  '
$ws.Range('E30').Style = 'Normal'

# Row 31
$ws.Range('D31').Value = '''4.080'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '''  +3.44%  '
$ws.Range('E31').Style = 'Normal'

# Row 32
$ws.Range('D32').Value = '''0.8952'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '''  +9.69%  '
$ws.Range('E32').Style = 'Normal'

# Row 33
$ws.Range('D33').Value = '''5.349'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '''  +8.04%  '
$ws.Range('E33').Style = 'Normal'

# Row 34
$ws.Range('D34').Value = '''0.08151'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '''  +2.73%  '
$ws.Range('E34').Style = 'Normal'

# Row 35
$ws.Range('D35').Value = '''1.594'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '''  +3.64%  '
$ws.Range('E35').Style = 'Normal'

# Row 36
$ws.Range('D36').Value = '''5.094'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '''  +8.01%  '
$ws.Range('E36').Style = 'Normal'

# Row 37
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').Value = '''1.232'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '''  +0.15%  '
$ws.Range('E37').Style = 'Normal'

# Row 38
$ws.Range('B38').Value = 'Aptos'
$ws.Range('C38').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D38').Value = '''11.45'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '''  +9.93%  '
$ws.Range('E38').Style = 'Normal'

# Row 39
$ws.Range('D39').Value = '''0.06007'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '''  +3.78%  '
$ws.Range('E39').Style = 'Normal'

# Row 40
$ws.Range('D40').Value = '''0.02161'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '''  +6.08%  '
$ws.Range('E40').Style = 'Normal'

# Row 41
$ws.Range('D41').Value = '''8.046'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '''  +5.92%  '
$ws.Range('E41').Style = 'Normal'

# Row 42
$ws.Range('D42').Value = '''0.1976'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '''  +5.34%  '
$ws.Range('E42').Style = 'Normal'

# Row 43
$ws.Range('D43').Value = '''0.9799'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '''  +2.28%  '
$ws.Range('E43').Style = 'Normal'

# Row 44
$ws.Range('D44').Value = '''0.5703'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '''  +8.20%  '
$ws.Range('E44').Style = 'Normal'

# Row 45
$ws.Range('D45').Value = '''12.87'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '''  +6.87%  '
$ws.Range('E45').Style = 'Normal'

# Row 46
$ws.Range('D46').Value = '''3.599'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '''  +2.46%  '
$ws.Range('E46').Style = 'Normal'

# Row 47
$ws.Range('D47').Value = '''0.5579'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '''  +7.54%  '
$ws.Range('E47').Style = 'Normal'

# Row 48
$ws.Range('D48').Value = '''123.93'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '''  +5.34%  '
$ws.Range('E48').Style = 'Normal'

# Row 49
$ws.Range('D49').Value = '''1.896'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '''  +6.10%  '
$ws.Range('E49').Style = 'Normal'

# Row 50
$ws.Range('D50').Value = '''0.06721'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '''  +4.21%  '
$ws.Range('E50').Style = 'Normal'

# Row 51
$ws.Range('D51').Value = '''71.44'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '''  +6.71%  '
$ws.Range('E51').Style = 'Normal'
